# Natmi following Dr Hou advice
# Rebuild the LR-pairs data: add "FAPs" as a second Sending cluster (in
# addition to the existing "sCs"), each paired with Wnt7b -> Fzd4 signalling
# to all three target clusters (ECs, FAPs, sCs), and refresh the computed
# statistics columns for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  A = "FAPs"; B = "Wnt7b"; C = "Fzd4"; D = "ECs";
       E = 1; F = 0.3333333333333333; G = 0.072919; H = 0.218757;
       I = 0.1016383815134179; J = 0.1016383815134179; K = 3; L = 1;
       M = 31.38723566666667; N = 94.16170700000001; O = 0.5539598599114094; P = 0.5539598599114095;
       Q = 2.288725837577667; R = 20.598532538199; S = 0.05630358358479538; T = 0.05630358358479539 },

    @{ Row = 3;  A = "FAPs"; B = "Wnt7b"; C = "Fzd4"; D = "FAPs";
       E = 1; F = 0.3333333333333333; G = 0.072919; H = 0.218757;
       I = 0.1016383815134179; J = 0.1016383815134179; K = 3; L = 1;
       M = 17.55525033333333; N = 52.665751; O = 0.3098362697066353; P = 0.3098362697066353;
       Q = 1.280111299056333; R = 11.521001691507; S = 0.03149125698713726; T = 0.03149125698713726 },

    @{ Row = 4;  A = "FAPs"; B = "Wnt7b"; C = "Fzd4"; D = "sCs";
       E = 1; F = 0.3333333333333333; G = 0.072919; H = 0.218757;
       I = 0.1016383815134179; J = 0.1016383815134179; K = 3; L = 1;
       M = 7.717279333333334; N = 23.151838; O = 0.1362038703819552; P = 0.1362038703819552;
       Q = 0.5627362917073333; R = 5.064626625366; S = 0.01384354094148529; T = 0.01384354094148529 },

    @{ Row = 5;  A = "sCs";  B = "Wnt7b"; C = "Fzd4"; D = "ECs";
       E = 3; F = 1; G = 0.6445166666666666; H = 1.93355;
       I = 0.8983616184865821; J = 0.898361618486582; K = 3; L = 1;
       M = 31.38723566666667; N = 94.16170700000001; O = 0.5539598599114094; P = 0.5539598599114095;
       Q = 20.22959650776111; R = 182.06636856985; S = 0.497656276326614; T = 0.4976562763266141 },

    @{ Row = 6;  A = "sCs";  B = "Wnt7b"; C = "Fzd4"; D = "FAPs";
       E = 3; F = 1; G = 0.6445166666666666; H = 1.93355;
       I = 0.8983616184865821; J = 0.898361618486582; K = 3; L = 1;
       M = 17.55525033333333; N = 52.665751; O = 0.3098362697066353; P = 0.3098362697066353;
       Q = 11.31465142733889; R = 101.83186284605; S = 0.2783450127194981; T = 0.278345012719498 },

    @{ Row = 7;  A = "sCs";  B = "Wnt7b"; C = "Fzd4"; D = "sCs";
       E = 3; F = 1; G = 0.6445166666666666; H = 1.93355;
       I = 0.8983616184865821; J = 0.898361618486582; K = 3; L = 1;
       M = 7.717279333333334; N = 23.151838; O = 0.1362038703819552; P = 0.1362038703819552;
       Q = 4.973915151655556; R = 44.7652363649; S = 0.1223603294404699; T = 0.1223603294404699 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $r.A
    $ws.Cells.Item($row, 2).Value  = $r.B
    $ws.Cells.Item($row, 3).Value  = $r.C
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
